$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Scanner" to "Pharmacology"
$ws.Name = "Pharmacology"

# Fix the Log Date column (C) formatting from MM/DD/YYYY to DD/MM/YYYY
# for all data rows (rows 2-40). Values are plain text, not real dates.
for ($r = 2; $r -le 40; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $old = [string]$cell.Text
    if ($old -match '^(\d{2})/(\d{2})/(\d{4})$') {
        $mm = $Matches[1]
        $dd = $Matches[2]
        $yyyy = $Matches[3]
        $cell.Value = "$dd/$mm/$yyyy"
    }
}
